$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new data row at row 66 (pushes the old "grand total" row 66
#     down to 67, and the old timestamp/footer row 67 down to 68) ---
$ws.Rows("66:66").Insert()

# Copy the formatting of the row above (row 65, a normal item row) onto the
# newly inserted blank row 66 so borders/fonts/number-formats match the rest
# of the item table.
$ws.Range("A65:Q65").Copy()
$ws.Range("A66:Q66").PasteSpecial(-4122)

# Match the row height used by the other item rows (row 65 uses 24.75,
# but this table alternates 24.75/25.5 - the inserted row takes 25.5,
# matching the size used elsewhere for this item slot).
$ws.Rows("66:66").RowHeight = 25.5

# --- Populate the new item row (item #60 - "مناديل FINE") ---
$ws.Range("A66").Value = 60
$ws.Range("C66").Value = "مناديل FINE"
$ws.Range("H66").Value = "29:0"

# L66, N66 and P66 hold numeric-looking text (stored as shared strings, not
# real numbers, in the source report) - force text storage with a "@" number
# format so Excel doesn't silently coerce them to numbers, then restore the
# original (copied) number format/style by re-pasting formats from row 65.
$ws.Range("L66").NumberFormat = "@"
$ws.Range("L66").Value = "0"
$ws.Range("L65").Copy()
$ws.Range("L66").PasteSpecial(-4122)

$ws.Range("N66").NumberFormat = "@"
$ws.Range("N66").Value = "30.00"
$ws.Range("N65").Copy()
$ws.Range("N66").PasteSpecial(-4122)

$ws.Range("P66").NumberFormat = "@"
$ws.Range("P66").Value = "30.0000"
$ws.Range("P65").Copy()
$ws.Range("P66").PasteSpecial(-4122)

$ws.Range("Q66").Value = "1:0"

# Re-create the merges for the new row (the blank inserted row has none).
$ws.Range("A66:B66").Merge()
$ws.Range("C66:G66").Merge()
$ws.Range("H66:K66").Merge()
$ws.Range("L66:M66").Merge()
$ws.Range("N66:O66").Merge()

# --- Update the grand-total row (old row 66, now row 67) with the new total ---
$ws.Range("P67").Value = 3312.7750000000001

# --- Update the footer timestamp (old row 67, now row 68) ---
$ws.Range("A68").Value = "Sunday, 7 September, 2025 7:12 PM"
